$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.246.24"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.536.25"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.26"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.21"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.537.42"
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("E11").Value = "  -5.30%  "
$ws.Range("E12").Value = "  +4.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.147.31"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.24"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.555.96"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.407.72"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.94"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.31"
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.37"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.686.48"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  +10.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +9.16%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.557.14"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("E36").Value = "  +7.60%  "
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.83"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.01"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.35"
$ws.Range("E43").Value = "  +17.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.93"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("E47").Value = "  +9.61%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.69"
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.457.88"
$ws.Range("E49").Value = "  +11.67%  "
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("E51").Value = "  +17.33%  "
